$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C, rows 2-11, from 45243 to 45244
for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 3).Value = 45244
}
